# Updated cryptos list values (Price & Volume(1h) columns) per source diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.041.33"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.832.19"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.64%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9994"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "241.57"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.37%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6554"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.37%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.000"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "44.56"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +5.86%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2932"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.72%  "
$ws.Range("E10").Value = "  -1.21%  "
$ws.Range("E11").Value = "  +0.53%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07680"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.53%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.841.55"
$ws.Range("D13").Style = "Normal"
$ws.Range("E14").Value = "  -0.58%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6668"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.16%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "81.52"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -5.32%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.116"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.09%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008676"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +4.30%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "29.051.07"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.37%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.089.73"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.73%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.43"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.68%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "224.28"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.95%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.000"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.08%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.132"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.83%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.000"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.03%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "157.51"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.84%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.496"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.22%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.1383"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.43%  "
$ws.Range("E29").Value = "  -0.53%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.507"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.26%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.109"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.68%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.201"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.97%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.007"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.50%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05359"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.92%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.841"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.01%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7415"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.10%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.158"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.92%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.651"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.97%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.296.60"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.37%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.755"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.92%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.01787"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.99%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.333"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +6.47%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9010"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.12%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.0000"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.22%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "103.09"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.26%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.989.00"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.78%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000124"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.85%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.07877"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.24%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "64.57"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.11%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.5140"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.43%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.741"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.88%  "
